$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 150 (shifts existing rows 150-193 down to 151-194)
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new data record
$ws.Cells.Item(150, 1).Value2 = 4
$ws.Cells.Item(150, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(150, 3).Value2 = "Los Lagos"
$ws.Cells.Item(150, 4).Value2 = 44588
$ws.Cells.Item(150, 5).Value2 = 10
$ws.Cells.Item(150, 6).Value2 = 100112017
$ws.Cells.Item(150, 7).Value2 = "Apio"
$ws.Cells.Item(150, 8).Value2 = "Americana (o)"
$ws.Cells.Item(150, 9).Value2 = "Primera"
$ws.Cells.Item(150, 10).Value2 = 25
$ws.Cells.Item(150, 11).Value2 = 12000
$ws.Cells.Item(150, 12).Value2 = 12000
$ws.Cells.Item(150, 13).Value2 = 12000
$ws.Cells.Item(150, 14).Value2 = "`$/docena de matas"
$ws.Cells.Item(150, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(150, 16).Value2 = 2000
$ws.Cells.Item(150, 17).Value2 = 6
$ws.Cells.Item(150, 18).Value2 = "Hortaliza"
